$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Three Season Tent" row (row 6) entirely - clears content and
# formatting so the row element is dropped from the sheet, and the two
# associated shared strings ("Three Season Tent" / "Comfortable shelter for
# two people.") are pruned since they become unused.
$ws.Range("A6:E6").Clear()

# Update the active selection to reflect where the user left off (A6).
$ws.Range("A6").Select()
